$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Part number renumbering (LT1... -> TGM1...), LBS/MC rev bump ---
$ws.Range("A3").Value = "TGM1A00-000"
$ws.Range("A4").Value = "TGM1A01-000"
$ws.Range("A5").Value = "TGM1A01-001"
$ws.Range("A6").Value = "TGM1A01-002"

# --- Enclosure -> MCU Unit rename ---
$ws.Range("B4").Value = "MCU Unit"

# --- New rows for Camera PCB / Camera Unit / Device Frame ---
$ws.Range("A7").Value = "TGM1A01-003"
$ws.Range("B7").Value = "Camera PCB"

$ws.Range("A17").Value = "TGM1A02-000"
$ws.Range("B17").Value = "Camera Unit"

$ws.Range("A25").Value = "TGM1A03-000"
$ws.Range("B25").Value = "Device Frame"

# --- Column A widens to fit the longer part numbers (was best-fit) ---
$ws.Columns.Item(1).ColumnWidth = 13.6

# --- Update selection to match authored state ---
$ws.Range("A26").Select()

$wb.Save()
